# Applies the commit "adding all the files when i submitted":
# a handful of adjacent same-formatted runs get merged back into single
# runs (normal Word re-save behaviour), the _GoBack bookmark moves from
# after "log_error is a two-line function," to inside "incomplete" in
# an earlier paragraph, an empty paragraph after the formula_to_graph
# bullet is removed, a new bullet point is appended to the "5. Notes"
# list, and the "Mohamed Hamza-vqwh89" / header runs are merged too.

$d = $word.ActiveDocument

$quoteL = [char]0x201C   # “
$quoteR = [char]0x201D   # ”

function Merge-Text($range, [string]$text) {
    # Find the given literal text in $range and "replace" it with the
    # exact same text; Word collapses the matched (possibly multi-run)
    # span into a single run sharing the first run's formatting.
    $null = $range.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

# --- word/document.xml -----------------------------------------------

# 1. "Alternatively, to install pygraphviz, install " + "brew from "
Merge-Text $d.Content "Alternatively, to install pygraphviz, install brew from "

# 2. Split "Or an incomplete one in the case..." after "Or an i" and
#    drop a _GoBack bookmark there (it gets relocated away from its old
#    spot automatically since a document only ever has one _GoBack).
$rng = $d.Content
$null = $rng.Find.Execute("Or an incomplete one in the case")
$splitPoint = $rng.Start + [string]"Or an i".Length
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 3. "A txt file: “grammar" + " " -> "A txt file: “grammar "
Merge-Text $d.Content ("A txt file: " + $quoteL + "grammar ")

# 4. "txt" + "”" + ", containing the grammar of this specific input file."
Merge-Text $d.Content ("txt" + $quoteR + ", containing the grammar of this specific input file.")

# 5. "userDefinedSyntax" + " parses the desired input file, "
Merge-Text $d.Content "userDefinedSyntax parses the desired input file, "

# 6. "log_error" + " is a two-line function," (old _GoBack bookmark here
#    already relocated by step 2 above)
Merge-Text $d.Content "log_error is a two-line function,"

# 7. "formula_to_graph" + " takes an instance of the formula class..."
Merge-Text $d.Content ("formula_to_graph takes an instance of the formula class and recursively expands on the object" + $quoteR.Replace($quoteR,[char]0x2019) + "s children, adding the relevant nodes to the networkX graph declared globally.")

# 7b. Remove the now-stray empty "ListParagraph" paragraph that used to
#     sit right after the formula_to_graph bullet (ind left=1440, no text).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "" -and $p.Range.ParagraphFormat.LeftIndent -eq 72) {
        $startText = $d.Range([Math]::Max(0, $p.Range.Start - 40), $p.Range.Start).Text
        if ($startText -like "*networkX graph declared globally.*") {
            $p.Range.Delete()
            break
        }
    }
}

# 8. Append the new "Note: ..." bullet at the very end of the document,
#    in the same numbered list (numId 8) as the preceding bullet, but
#    with a smaller font size.
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$apost = [char]0x2019
$newPara.Range.Text = "Note: the grammar states formula -> (formula), this is allowed as all parentheses are stripped in my implementation, I did not realise the invalidity that it would make in your FOL grammar, and assumed that since it makes sense semantically it would take too much re-working of the formula class to find redundant braces for only certain cases."
$newPara.Range.Font.Size = 7.5
$newPara.Range.Font.SizeBi = 7.5

# --- word/footer2.xml --------------------------------------------------

$sec = $d.Sections.Item(1)
$primaryFooter = $sec.Footers.Item(1)
Merge-Text $primaryFooter.Range "Mohamed Hamza-vqwh89"

# --- word/header1.xml ---------------------------------------------------

$primaryHeader = $sec.Headers.Item(1)
Merge-Text $primaryHeader.Range "COMP2211 Networks and Systems"
Merge-Text $primaryHeader.Range "Compiler Design Assignment Documentation"

Write-Output "done"
